# msPRR_01 correct file key
# Rename every shared-string cell value of the form "msPRR-01-<N>" to
# "msPRR_01_<N>" (dash -> underscore) throughout the worksheet, and
# restore the sheet view (zoom + active selection) recorded by the
# author at commit time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$firstRow = $used.Row
$firstCol = $used.Column

for ($r = $firstRow; $r -lt ($firstRow + $rowCount); $r++) {
    for ($c = $firstCol; $c -lt ($firstCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string] -and $v.StartsWith("msPRR-01-")) {
            $cell.Value2 = $v.Replace("msPRR-01-", "msPRR_01_")
        }
    }
}

# Restore the view state: zoomed to 215%, scrolled so column K is at the
# left edge, with O3 as the active selection.
$win = $excel.ActiveWindow
$ws.Range("O3").Select()
$win.ScrollColumn = 11
$win.ScrollRow = 1
$win.Zoom = 215
